$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1:H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for column I and J, rows 2..36
$data = @(
    @(6, 7),
    @(7, 9),
    @(9, 9),
    @(9, 9),
    @(4, 5),
    @(2, 4),
    @(6, 6),
    @(3, 5),
    @(9, 9),
    @(6, 8),
    @(7, 7),
    @(8, 8),
    @(1, 4),
    @(1, 6),
    @(1, 4),
    @(1, 5),
    @(1, 3),
    @(1, 6),
    @(1, 3),
    @(1, 4),
    @(1, 5),
    @(1, 1),
    @(1, 5),
    @(1, 5),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 3),
    @(1, 4),
    @(1, 5),
    @(1, 4),
    @(1, 4),
    @(1, 4),
    @(1, 2)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
